# تعديل يدوي في شيت Card21
# Move the "عطل تلامس متكرر / T.CON" incident text (and its long note) from row 29
# down into row 30, clear the now-vacated cells in row 29, and clear O31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Text that needs to move from row 29 to row 30 ---
$faultText = "عطل تلامس متكرر`nT.CON"
$noteText  = "تم تغيير الجرائد الاماميه عند طن 979`nو معيارية الجرائد الاماميه علي 15 ساو`nو اليروفيل الامامي علي 60 ساو`nو السكينه الاماميه علي 11 ساو `nو اليروفيل الخلفي علي 60/60/70 ساو`nوالسكينه الخلفيه  20ساو"

# --- Row 30 gets the moved text in M30/N30; L30 and O30 stay as they were ---
$ws.Range("M30").Value = $faultText
$ws.Range("N30").Value = $noteText

# --- Row 29: L29, M29, N29, O29 are cleared out (content moved away / removed) ---
$ws.Range("L29").ClearContents()
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("O29").ClearContents()

# --- Row 31: O31 is cleared ---
$ws.Range("O31").ClearContents()
